# Apply "Arch.xlsx" / effec_rec_field sheet updates:
#  - extend the B->G blocks for "receptive field" tables down to row 56
#  - switch the G43:G47 formulas to the new F*(C-1)+1 form
#  - fill in row 48 (previously blank) as a continuation of the first table
#  - rescale the second table (rows 51-56) for a 200x200x50 scene,
#    adding a new row 56

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("effec_rec_field")

# --- Rows 43-47: update the G formula (keep the same cached values) ---
$ws.Range("G43").Formula = "=(B43-1)*E43-2*D43+F43*(C43-1)+1"
$ws.Range("G44").Formula = "=(B44-1)*E44-2*D44+F44*(C44-1)+1"
$ws.Range("G45").Formula = "=(B45-1)*E45-2*D45+F45*(C45-1)+1"
$ws.Range("G46").Formula = "=(B46-1)*E46-2*D46+F46*(C46-1)+1"
$ws.Range("G47").Formula = "=(B47-1)*E47-2*D47+F47*(C47-1)+1"

# --- Row 48: was blank, now a continuation row of the first table ---
$ws.Range("B48").Formula = "=G47"
$ws.Range("C48").Value = 4
$ws.Range("D48").Value = 1
$ws.Range("E48").Value = 2
$ws.Range("F48").Value = 1
$ws.Range("G48").Formula = "=(B48-1)*E48-2*D48+F48*(C48-1)+1"
$ws.Range("B48:G48").HorizontalAlignment = -4108

# --- Second table (rows 51-55): rescale for a 200x200x50 scene ---
$ws.Range("B51").Value = 200

$ws.Range("E52").Value = 2
$ws.Range("E53").Value = 2

# Row 55 gains the values that used to belong to the (not yet existing) row 56
$ws.Range("E55").Value = 3

# --- New row 56, continuing the second table ---
$ws.Range("B56").Formula = "=G55"
$ws.Range("C56").Value = 4
$ws.Range("D56").Value = 1
$ws.Range("E56").Value = 3
$ws.Range("F56").Value = 1
$ws.Range("G56").Formula = "=ROUNDDOWN((B56+2*D56-F56*(C56-1)-1)/E56 + 1,0)"
$ws.Range("B56:G56").HorizontalAlignment = -4108

# --- Selection, matching the saved view state in the edited workbook ---
$ws.Range("E37").Select()
